$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 94 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
